$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Host "NOT FOUND: $find"
    }
}

# 1. Objective paragraph wording change
Replace-Text "working towards real goals" "working to create beautiful interfaces which are a joy to use"

# 2. Skills: react-native, electron or ionic -> React Native, Electron or Ionic
Replace-Text "react-native, electron or ionic" "React Native, Electron or Ionic"

# 3. Skills: leaning new ones -> learning new ones
Replace-Text "leaning new ones" "learning new ones"

# 4. Secret Lab dates: (Nov 18 - Present) -> (Nov 18 - Jun 2020)
Replace-Text "(Nov 18 - Present)" "(Nov 18 - Jun 2020)"

# 5. Title: Senior React-Native Developer -> Senior React Native Developer
Replace-Text "Senior React-Native Developer" "Senior React Native Developer"

# 6. Bullet: using react-native for the front-end -> using React Native for the front-end
Replace-Text "using react-native for the front-end" "using React Native for the front-end"

# 7a. Teaching Assistant bullet 1 wording
Replace-Text "We focus primarily on the MERN stack. We covered topics" "We focused primarily on the MERN stack and covered topics"

# 7b. Teaching Assistant bullet 2 typo fix
Replace-Text "Through One-on-on and group mentoring" "Through One-on-one and group mentoring"

# 8. Thomson Reuters bullet: Jquery -> jQuery (Backbone line)
Replace-Text "Backbone, Jquery and SCSS" "Backbone, jQuery and SCSS"

# 9. ScanCafe bullet: Zend, and Jquery. -> Zend, and jQuery.
Replace-Text "Zend, and Jquery." "Zend, and jQuery."
